# Ran code for averaged intensities on spiral schemes.
#
# The previously-last scheme row ("Gaussian-Quadrature", old row 16) moves up
# to row 10. Three new "Spiral" rows are inserted right after it (rows
# 11-13), and the schemes that used to occupy rows 10-15
# ("NoRotation-tilt60deg" ... "HexGrid-60degTilt5degRes") shift down to rows
# 14-19, renumbered accordingly in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- make room: insert 3 blank rows before the current row 11 ---
$ws.Rows("11:13").Insert()

# copy the column-A style (bold/border) down onto the newly inserted rows
$ws.Range("A10").Copy()
$ws.Range("A11:A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P")

function Set-SchemeRow {
    param($row, $idx, $label, $vals)
    $ws.Range("A$row").Value = $idx
    $ws.Range("B$row").Value = $label
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

# --- row 10: Gaussian-Quadrature (relocated from the old row 16) ---
Set-SchemeRow 10 8 "Gaussian-Quadrature" @(1.006681252277507,0.9668042872944713,0.9982645422461321,0.9918008236398924,1.006681252277507,0.9668042872944713,1.002759702876757,0.9879085223676277,0.998656679882176,0.9761799342428142,1.006681252277507,0.9825344147703017,0.9908877263645007,0.9911319681034223)

# --- row 11: Spiral-90deg-10rot-5space (new) ---
Set-SchemeRow 11 9 "Spiral-90deg-10rot-5space" @(1.037252785942101,0.9257705356791749,1.007978139672298,0.9795215192644341,1.037252785942101,0.9257705356791749,1.015973468772326,0.9794944047942302,1.007946382976398,0.9491220813218516,1.037252785942101,0.9668743376757367,0.9876307451395021,0.9878824148028518)

# --- row 12: Spiral-90deg-15rot-5space (new) ---
Set-SchemeRow 12 10 "Spiral-90deg-15rot-5space" @(1.037011990648133,0.9260745748039266,1.007909046333152,0.9795710240738237,1.037011990648133,0.9260745748039266,1.015853936067461,0.9795534416124833,1.007857431248952,0.9493100857137541,1.037011990648133,0.9669918105685391,0.9876416589647588,0.9878926913127106)

# --- row 13: Spiral-90deg-10rot-3space (new) ---
Set-SchemeRow 13 11 "Spiral-90deg-10rot-3space" @(1.037130940620655,0.9258562276167617,1.0079537762545,0.9795288886891909,1.037130940620655,0.9258562276167617,1.01594740929233,0.9795134664762274,1.007926333052678,0.9491646545956037,1.037130940620655,0.9669050019356307,0.9876174582952768,0.9878777120747433)

# --- row 14: NoRotation-tilt60deg (relocated from the old row 10) ---
Set-SchemeRow 14 12 "NoRotation-tilt60deg" @(1.053571999999998,0.9164279999999991,1.016883999999999,0.976664,1.053571999999998,0.9164279999999991,1.017927999999999,0.977700000000002,1.010424000000001,0.9489400000000009,1.053571999999998,0.9666559999999993,0.9908869999999991,0.9898175)

# --- row 15: Rotation-NoTilt (relocated from the old row 11) ---
Set-SchemeRow 15 13 "Rotation-NoTilt" @(1.09,0.88,1.03,0.97,1.09,0.88,1.03,0.97,1.02,0.93,1.09,0.9550000000000001,0.9924999999999999,0.99)

# --- row 16: Rotation-60detTilt (relocated from the old row 12) ---
Set-SchemeRow 16 14 "Rotation-60detTilt" @(1.047306194943997,0.9276161642496022,1.014296358502399,0.9797005072384041,1.047306194943997,0.9276161642496022,1.014398557183996,0.9798552248320014,1.008841585254397,0.9565291417599987,1.047306194943997,0.9709562613760006,0.9922298062336005,0.9910679667455995)

# --- row 17: HexGrid-90degTilt5degRes (relocated from the old row 13) ---
Set-SchemeRow 17 15 "HexGrid-90degTilt5degRes" @(0.9930148727799156,0.992850651667461,0.9930417773476204,0.9931438929443994,0.9930148727799156,0.992850651667461,0.9927537972816609,0.9931960738164011,0.9934195112590208,0.9929055828988569,0.9930148727799156,0.9929462145075407,0.9930127986848492,0.993040769999417)

# --- row 18: HexGrid-90degTilt22p5degRes (relocated from the old row 14) ---
Set-SchemeRow 18 16 "HexGrid-90degTilt22p5degRes" @(0.9905640992280428,0.9954270959934313,0.9938167432511223,0.9941891180771878,0.9905640992280428,0.9954270959934313,0.9910132575558535,0.993428950588199,0.9923700646700258,0.9945724463915644,0.9905640992280428,0.9946219196222768,0.993499264137446,0.9931727219694284)

# --- row 19: HexGrid-60degTilt5degRes (relocated from the old row 15) ---
Set-SchemeRow 19 17 "HexGrid-60degTilt5degRes" @(0.9848088518153224,1.004004455682517,0.9907288351983266,0.9955773762197788,0.9848088518153224,1.004004455682517,0.9883989982630875,0.9955606735335814,0.9906183857939637,1.000638715397396,0.9848088518153224,0.9973666454404215,0.9937798797289861,0.9937920364879966)

Write-Host "Done applying spiral-scheme averaged intensities edit."
